$d = $word.ActiveDocument

function FindReplace([string]$find, [string]$replace) {
    $ok = $d.Content.Find.Execute($find, $false, $true, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $find"
    }
    return $ok
}

function FindParagraphIndex([string]$exactText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $exactText) {
            return $i
        }
    }
    return -1
}

# 1. "Simple Traffic Sim V1" -> "Simple Traffic Sim V3"
FindReplace "Simple Traffic Sim V1" "Simple Traffic Sim V3"

# 2. Append a new sentence about the "egui" user interface at the end of the
#    Overview section's first paragraph.
FindReplace "Objects can only be spawned if the grid is not occupied at that location. " `
    "Objects can only be spawned if the grid is not occupied at that location. There is a user interface implemented with “egui” that provides buttons, controls, and stats about the game. "

# 3. Vehicles section: remove "which should be replaced with a better algorithm"
#    and merge the following paragraph in, adding a new sentence about lane choice.
FindReplace "do a simple DFS for the path, which should be replaced with a better algorithm. If the path is edited" `
    "do a simple DFS for the path. If the path is edited"

FindReplace "relevant graph components are changed.^pVehicles drive freely and follow the desired lane on the current road. Vehicles must turn manually" `
    "relevant graph components are changed. Vehicles drive freely and follow the desired lane on the current road. The desired lane is based on the upcoming turn direction in their search path. Vehicles must turn manually"

# 4. Insert a new "User Interface" section (heading + paragraph) right before
#    the "Saving" heading.
$savingIdx = FindParagraphIndex "Saving"
if ($savingIdx -lt 0) {
    Write-Output "WARNING: could not find Saving heading paragraph"
} else {
    $savingPara = $d.Paragraphs.Item($savingIdx)
    $savingPara.Range.InsertParagraphBefore()
    $savingPara.Range.InsertParagraphBefore()

    $headingPara = $d.Paragraphs.Item($savingIdx)
    $bodyPara = $d.Paragraphs.Item($savingIdx + 1)

    $headingPara.Range.Text = "User Interface"
    $headingPara.Range.Bold = 1

    $bodyPara.Range.Text = "The game integrates the “egui” system using the “bevy_egui” crate. I used this to add a very basic user interface for switching between tools, describing controls, and showing stats about the game state. Although the immediate mode GUI is not great for a game application, it is a good start for getting some basic info on screen and can be used for development tools. Additionally, the bevy engine’s native GUI system is not particularly powerful or concise, so I am avoiding it."
    $bodyPara.Range.Bold = 0
}

# 5. "How to play the demo" paragraph updates.
FindReplace "Hold [L] to spawn many vehicles." "Vehicles will begin spawning automatically."

FindReplace "Observe that vehicles detect each other and slow down to prevent collisions. Press [F5]" `
    "Observe that vehicles detect each other and slow down to prevent collisions. The vehicles should also change lanes based on what turn they need to make. Press [F5]"

Write-Output "done"
